$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Col=4; Val="'47.315.59"},
    @{Row=2; Col=5; Val="  +2.35%  "},
    @{Row=3; Col=4; Val="'2.501.93"},
    @{Row=3; Col=5; Val="  +2.11%  "},
    @{Row=4; Col=5; Val="  +0.09%  "},
    @{Row=5; Col=4; Val="'323.52"},
    @{Row=5; Col=5; Val="  +0.51%  "},
    @{Row=6; Col=4; Val="'108.77"},
    @{Row=6; Col=5; Val="  +3.60%  "},
    @{Row=7; Col=5; Val="  +1.32%  "},
    @{Row=8; Col=4; Val="'1.00"},
    @{Row=8; Col=5; Val="  -0.02%  "},
    @{Row=9; Col=4; Val="'0.536"},
    @{Row=9; Col=5; Val="  -0.14%  "},
    @{Row=10; Col=4; Val="'39.06"},
    @{Row=10; Col=5; Val="  +8.47%  "},
    @{Row=11; Col=4; Val="'0.0813"},
    @{Row=11; Col=5; Val="  +0.77%  "},
    @{Row=12; Col=5; Val="  +0.76%  "},
    @{Row=13; Col=4; Val="'18.37"},
    @{Row=13; Col=5; Val="  +0.43%  "},
    @{Row=14; Col=4; Val="'7.18"},
    @{Row=14; Col=5; Val="  +1.48%  "},
    @{Row=15; Col=4; Val="'2.894.74"},
    @{Row=15; Col=5; Val="  +2.13%  "},
    @{Row=16; Col=4; Val="'2.511.49"},
    @{Row=16; Col=5; Val="  +1.89%  "},
    @{Row=17; Col=4; Val="'0.854"},
    @{Row=17; Col=5; Val="  +1.42%  "},
    @{Row=18; Col=4; Val="'47.242.93"},
    @{Row=18; Col=5; Val="  +2.46%  "},
    @{Row=19; Col=4; Val="'12.81"},
    @{Row=19; Col=5; Val="  +1.43%  "},
    @{Row=20; Col=4; Val="'6.64"},
    @{Row=20; Col=5; Val="  +3.48%  "},
    @{Row=21; Col=4; Val="'0.0₃0942"},
    @{Row=21; Col=5; Val="  +0.63%  "},
    @{Row=22; Col=4; Val="'2.70"},
    @{Row=22; Col=5; Val="  +12.74%  "},
    @{Row=23; Col=4; Val="'70.50"},
    @{Row=24; Col=4; Val="'247.82"},
    @{Row=24; Col=5; Val="  +0.17%  "},
    @{Row=25; Col=4; Val="'2.60"},
    @{Row=25; Col=5; Val="  +3.02%  "},
    @{Row=26; Col=4; Val="'26.06"},
    @{Row=26; Col=5; Val="  +0.67%  "},
    @{Row=27; Col=4; Val="'1.00"},
    @{Row=27; Col=5; Val="  -0.02%  "},
    @{Row=28; Col=5; Val="  +0.33%  "},
    @{Row=29; Col=4; Val="'10.07"},
    @{Row=29; Col=5; Val="  +3.90%  "},
    @{Row=30; Col=4; Val="'35.29"},
    @{Row=30; Col=5; Val="  +2.92%  "},
    @{Row=31; Col=4; Val="'0.138"},
    @{Row=31; Col=5; Val="  +6.81%  "},
    @{Row=32; Col=4; Val="'49.81"},
    @{Row=32; Col=5; Val="  +0.85%  "},
    @{Row=33; Col=4; Val="'19.99"},
    @{Row=33; Col=5; Val="  +0.91%  "},
    @{Row=34; Col=4; Val="'5.43"},
    @{Row=34; Col=5; Val="  +1.68%  "},
    @{Row=35; Col=4; Val="'0.0790"},
    @{Row=35; Col=5; Val="  +3.34%  "},
    @{Row=37; Col=4; Val="'1.98"},
    @{Row=37; Col=5; Val="  +4.61%  "},
    @{Row=38; Col=4; Val="'4.69"},
    @{Row=38; Col=5; Val="  +3.43%  "},
    @{Row=39; Col=4; Val="'2.99"},
    @{Row=39; Col=5; Val="  +1.32%  "},
    @{Row=40; Col=5; Val="  +1.07%  "},
    @{Row=41; Col=5; Val="  +0.68%  "},
    @{Row=42; Col=4; Val="'121.03"},
    @{Row=42; Col=5; Val="  -5.53%  "},
    @{Row=43; Col=4; Val="'21.41"},
    @{Row=43; Col=5; Val="  +2.71%  "},
    @{Row=44; Col=4; Val="'0.0299"},
    @{Row=45; Col=4; Val="'1.991.71"},
    @{Row=45; Col=5; Val="  +1.08%  "},
    @{Row=46; Col=4; Val="'3.07"},
    @{Row=46; Col=5; Val="  +3.22%  "},
    @{Row=47; Col=5; Val="  -1.48%  "},
    @{Row=48; Col=2; Val="FraxShare"},
    @{Row=48; Col=3; Val="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"},
    @{Row=48; Col=4; Val="'9.10"},
    @{Row=48; Col=5; Val="  +0.20%  "},
    @{Row=49; Col=2; Val="Stacks"},
    @{Row=49; Col=3; Val="https://coinranking.com/coin/mMPrMcB7+stacks-stx"},
    @{Row=49; Col=4; Val="'1.78"},
    @{Row=49; Col=5; Val="  -4.28%  "},
    @{Row=50; Col=4; Val="'5.21"},
    @{Row=50; Col=5; Val="  +3.15%  "},
    @{Row=51; Col=4; Val="'56.64"},
    @{Row=51; Col=5; Val="  +4.01%  "}
)

foreach ($item in $updates) {
    $ws.Cells.Item($item.Row, $item.Col).Value = $item.Val
}

Write-Output "Applied $($updates.Count) cell updates"